$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "time_taken" header in F1, matching the style used by the other
# header cells (bold, centered, bordered) by copying the format from E1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Populate the time_taken column for each data row with plain (unstyled)
# timestamp text values.
$ws.Range("F2").Value = "2021-10-05 13:42:27.783732"
$ws.Range("F3").Value = "2021-10-05 13:42:27.783742"
$ws.Range("F4").Value = "2021-10-05 13:42:27.783746"
$ws.Range("F5").Value = "2021-10-05 13:42:27.783748"
$ws.Range("F6").Value = "2021-10-05 13:42:27.783751"
